$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.445.98'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '3.976.14'
$ws.Range('E3').Value = '  -2.03%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '610.79'
$ws.Range('E5').Value = '  +7.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.87'
$ws.Range('E6').Value = '  +11.42%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.680'
$ws.Range('E7').Value = '  -2.55%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.780'
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.185'
$ws.Range('E10').Value = '  +6.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.21'
$ws.Range('E11').Value = '  +3.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000334'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.27'
$ws.Range('E13').Value = '  +1.29%  '
$ws.Range('D14').Value = '4.607.46'
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('D15').Value = '3.962.55'
$ws.Range('E15').Value = '  -2.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.20'
$ws.Range('E16').Value = '  -2.42%  '
$ws.Range('E17').Value = '  +1.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.51'
$ws.Range('E18').Value = '  -1.68%  '
$ws.Range('D19').Value = '73.265.57'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '467.89'
$ws.Range('E21').Value = '  +4.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.83'
$ws.Range('E22').Value = '  +8.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '96.66'
$ws.Range('E23').Value = '  -1.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.42'
$ws.Range('E24').Value = '  -4.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.19'
$ws.Range('E25').Value = '  -4.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.22'
$ws.Range('E26').Value = '  -3.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.02'
$ws.Range('E27').Value = '  -2.84%  '
$ws.Range('E28').Value = '  +0.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.51'
$ws.Range('E29').Value = '  -5.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.29'
$ws.Range('E30').Value = '  -2.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.88'
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.91'
$ws.Range('E32').Value = '  +2.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0000104'
$ws.Range('E33').Value = '  +14.97%  '
$ws.Range('E34').Value = '  -4.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '47.90'
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '70.65'
$ws.Range('E36').Value = '  +3.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '648.47'
$ws.Range('E37').Value = '  -5.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.430'
$ws.Range('E38').Value = '  -4.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.39'
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.16'
$ws.Range('E44').Value = '  +36.03%  '
$ws.Range('B45').Value = 'THORChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.45'
$ws.Range('E45').Value = '  -6.28%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.14'
$ws.Range('E46').Value = '  -5.99%  '
$ws.Range('E47').Value = '  -2.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000296'
$ws.Range('E48').Value = '  +5.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.42'
$ws.Range('E49').Value = '  +3.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.56'
$ws.Range('E50').Value = '  -4.87%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.827.25'
$ws.Range('E51').Value = '  +1.79%  '
